$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    17662.27096434525,
    17332.44617459095,
    17020.93522754256,
    16833.32822805317,
    16538.7342770511,
    16251.78958823426,
    17048.68430782569,
    16867.07296985057,
    16575.0902168118,
    16763.74460718088,
    16351.94345054924,
    16314.91758482173,
    16162.77515283737,
    15965.8305228978,
    15849.80456782925,
    15817.51626305694,
    15701.78034877345,
    15601.36284657243,
    15528.12566353889,
    15433.95384278634,
    15294.55460663398,
    15158.90057301816,
    15030.45763530442,
    14857.28982172285,
    14845.0408508628,
    14753.63493998294,
    14662.02616363277,
    14579.83724974177,
    14514.50781988384,
    14446.00062549772,
    14365.10377217448,
    14294.48538702716,
    14237.04611460397,
    14230.64640748517,
    14220.5246005539,
    14151.2143321999,
    14155.34556206754,
    14093.76799487183,
    14058.43545997221,
    13945.55298608913,
    13970.3714367555,
    13939.64634520983,
    13895.86956098587,
    13876.01666224103,
    13838.9260151458,
    13748.83502408908,
    13757.16632557489,
    13723.26536673252,
    13683.6258856331,
    13630.64496857635,
    13594.76007637084,
    13540.00493931051,
    13499.49915695492,
    13448.66614774409,
    13390.78340461597,
    13386.71993698526,
    13340.24435578429,
    13299.53133240902
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
